$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct "Irrazabal" -> "Irrázabal" (row 10 / id_usuario 8, segundo_apellido)
$ws.Range("E10").Value = "Irrázabal"

# Correct "Premolo" -> "Prémolo" (row 18 / id_usuario 16, primer_apellido)
$ws.Range("D18").Value = "Prémolo"

# Correct "Gutierrez" -> "Gutiérrez" (row 24 / id_usuario 22, primer_apellido)
$ws.Range("D24").Value = "Gutiérrez"

# Correct admin user's password hash (row 3 / id_usuario 1)
$ws.Range("C3").Value = "c16fd958b85a1c94d872c219ea06ce8e80223239b1fcefb92ad978445ef095507244be44caae1d766e277b072c184cb3ffe4d0610716e989b2fe5a7c97bf3144"

$wb.Application.Calculate()
